$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.885.17'
$ws.Range("E2").Value = '  +3.90%  '
$ws.Range("D3").Value = '2.525.40'
$ws.Range("E3").Value = '  +2.15%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.38'
$ws.Range("E5").Value = '  +3.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.99'
$ws.Range("E6").Value = '  +5.96%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +3.50%  '
$ws.Range("D9").Value = '2.524.56'
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("E10").Value = '  +4.72%  '
$ws.Range("E11").Value = '  +3.42%  '
$ws.Range("E12").Value = '  +2.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.84'
$ws.Range("E14").Value = '  +3.22%  '
$ws.Range("D15").Value = '2.982.89'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("E16").Value = '  +4.08%  '
$ws.Range("D17").Value = '67.715.34'
$ws.Range("E17").Value = '  +3.45%  '
$ws.Range("D18").Value = '2.512.09'
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.99'
$ws.Range("E19").Value = '  +6.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.45'
$ws.Range("E20").Value = '  +3.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '360.01'
$ws.Range("E21").Value = '  +5.74%  '
$ws.Range("E22").Value = '  +1.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.66'
$ws.Range("E23").Value = '  +3.96%  '
$ws.Range("E24").Value = '  +2.43%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.18'
$ws.Range("E26").Value = '  +4.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.26'
$ws.Range("E27").Value = '  +5.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("D29").Value = '2.644.03'
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").Value = '0.0₃0992'
$ws.Range("E30").Value = '  +4.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '543.43'
$ws.Range("E31").Value = '  +5.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.28'
$ws.Range("E32").Value = '  +4.38%  '
$ws.Range("E33").Value = '  +3.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  +4.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.129'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.47'
$ws.Range("E37").Value = '  +2.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '155.15'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.82'
$ws.Range("E39").Value = '  +2.94%  '
$ws.Range("E40").Value = '  +2.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.357'
$ws.Range("E41").Value = '  +2.38%  '
$ws.Range("E42").Value = '  +3.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.18'
$ws.Range("E43").Value = '  +4.52%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  +5.89%  '
$ws.Range("E46").Value = '  +2.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '146.50'
$ws.Range("E47").Value = '  +1.09%  '
$ws.Range("E48").Value = '  +5.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.73'
$ws.Range("E49").Value = '  +2.73%  '
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("E51").Value = '  +1.91%  '
